$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 837882.9399999999
$ws.Range("J40").Value = 7000
$ws.Range("L40").Value = 7000
$ws.Range("N40").Value = -7350
$ws.Range("H80").Value = 4405.8125
$ws.Range("J80").Value = 4899.1
$ws.Range("L80").Value = 14697.3
$ws.Range("N80").Value = -16693.3
$ws.Range("H83").Value = 4405.8125
$ws.Range("J83").Value = 4899.1
$ws.Range("L83").Value = 44091.9
$ws.Range("N83").Value = -54075.9
$ws.Range("H112").Value = 6330972
$ws.Range("J112").Value = 6330972
$ws.Range("L112").Value = 18992916
$ws.Range("N112").Value = -18995132
$ws.Range("H116").Value = 4032.6667
$ws.Range("I116").Value = 4064.3333
$ws.Range("J116").Value = 4001
$ws.Range("K116").Value = 4064.3333
$ws.Range("L116").Value = 4001
$ws.Range("M116").Value = -622.3332999999998
$ws.Range("N116").Value = -10885
$ws.Range("H132").Value = 14286820
$ws.Range("I132").Value = 14926405
$ws.Range("K132").Value = 44779215
$ws.Range("M132").Value = -44776685
$ws.Range("H137").Value = 16970.818
$ws.Range("I137").Value = 5248
$ws.Range("K137").Value = 15744
$ws.Range("M137").Value = -13194
$ws.Range("H138").Value = 3808.641
$ws.Range("I138").Value = 859.0714
$ws.Range("K138").Value = 2577.2142
$ws.Range("M138").Value = 2562.7858
$ws.Range("H141").Value = 969.78845
$ws.Range("I141").Value = 884.38
$ws.Range("K141").Value = 2653.14
$ws.Range("M141").Value = 2526.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5997.4
$ws.Range("H31").Value = 131409.2
$ws.Range("I31").Value = 21368.75
$ws.Range("J31").Value = 571571
$ws.Range("K31").Value = 21368.75
$ws.Range("L31").Value = 571571
$ws.Range("M31").Value = -21074.75
$ws.Range("N31").Value = -572159
$ws.Range("H32").Value = 4733.5527
$ws.Range("I32").Value = 3742.0303
$ws.Range("J32").Value = 11277.6
$ws.Range("K32").Value = 3742.0303
$ws.Range("L32").Value = 11277.6
$ws.Range("M32").Value = -3455.0303
$ws.Range("N32").Value = -11851.6
$ws.Range("H45").Value = 2783
$ws.Range("I45").Value = 1900
$ws.Range("J45").Value = 3003.75
$ws.Range("K45").Value = 1900
$ws.Range("L45").Value = 3003.75
$ws.Range("M45").Value = -1523
$ws.Range("N45").Value = -3757.75
$ws.Range("H74").Value = 43528150
$ws.Range("I74").Value = 50057084
$ws.Range("K74").Value = 50057084
$ws.Range("M74").Value = -50056210
$ws.Range("H77").Value = 43528150
$ws.Range("I77").Value = 50057084
$ws.Range("K77").Value = 250285420
$ws.Range("M77").Value = -250281052
$ws.Range("H110").Value = 2064.5454
$ws.Range("I110").Value = 1401.375
$ws.Range("K110").Value = 1401.375
$ws.Range("M110").Value = 643.625
$ws.Range("H116").Value = 5997.4
$ws.Range("H132").Value = 17591564
$ws.Range("I132").Value = 6127.7144
$ws.Range("K132").Value = 18383.1432
$ws.Range("M132").Value = -15853.1432
$ws.Range("H138").Value = 199999
$ws.Range("J138").Value = 199999
$ws.Range("L138").Value = 199999
$ws.Range("N138").Value = -210279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5997.4
$ws.Range("H86").Value = 22278.334
$ws.Range("J86").Value = 60743.145
$ws.Range("L86").Value = 60743.145
$ws.Range("N86").Value = -62989.145
$ws.Range("H89").Value = 22278.334
$ws.Range("J89").Value = 60743.145
$ws.Range("L89").Value = 303715.725
$ws.Range("N89").Value = -314947.725
$ws.Range("H105").Value = 9505.5
$ws.Range("I105").Value = 15154.714
$ws.Range("J105").Value = 3856.2856
$ws.Range("K105").Value = 15154.714
$ws.Range("L105").Value = 3856.2856
$ws.Range("M105").Value = -13407.714
$ws.Range("N105").Value = -7350.2856
$ws.Range("H107").Value = 4049.1428
$ws.Range("I107").Value = 3669
$ws.Range("K107").Value = 3669
$ws.Range("M107").Value = -1749

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1415.909
$ws.Range("I16").Value = 1098.1666
$ws.Range("K16").Value = 1098.1666
$ws.Range("M16").Value = -811.1666
$ws.Range("H31").Value = 28575420
$ws.Range("I31").Value = 2464.2942
$ws.Range("K31").Value = 2464.2942
$ws.Range("M31").Value = -2169.2942
$ws.Range("H34").Value = 28575420
$ws.Range("I34").Value = 2464.2942
$ws.Range("K34").Value = 2464.2942
$ws.Range("M34").Value = -2262.2942
$ws.Range("H58").Value = 3413.15
$ws.Range("I58").Value = 2857.4
$ws.Range("J58").Value = 3968.9
$ws.Range("K58").Value = 2857.4
$ws.Range("L58").Value = 3968.9
$ws.Range("M58").Value = -2654.4
$ws.Range("N58").Value = -4374.9
$ws.Range("H113").Value = 1415.909
$ws.Range("I113").Value = 1098.1666
$ws.Range("K113").Value = 1098.1666
$ws.Range("M113").Value = 1071.8334
$ws.Range("H132").Value = 90935.69500000001
$ws.Range("I132").Value = 105488.38
$ws.Range("K132").Value = 316465.14
$ws.Range("M132").Value = -313935.14
$ws.Range("H136").Value = 3413.15
$ws.Range("I136").Value = 2857.4
$ws.Range("J136").Value = 3968.9
$ws.Range("K136").Value = 8572.200000000001
$ws.Range("L136").Value = 11906.7
$ws.Range("M136").Value = -6022.200000000001
$ws.Range("N136").Value = -17006.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 154.45454
$ws.Range("I11").Value = 161.9
$ws.Range("K11").Value = 485.7
$ws.Range("M11").Value = -345.7
$ws.Range("H26").Value = 493.57144
$ws.Range("I26").Value = 177.5
$ws.Range("K26").Value = 532.5
$ws.Range("M26").Value = -244.5
$ws.Range("H132").Value = 2474750.2
$ws.Range("I132").Value = 1398.5555
$ws.Range("J132").Value = 3711426
$ws.Range("K132").Value = 12586.9995
$ws.Range("L132").Value = 33402834
$ws.Range("M132").Value = -10056.9995
$ws.Range("N132").Value = -33407894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 8225.362999999999
$ws.Range("I99").Value = 8225.362999999999
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8225.362999999999
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -5979.362999999999
$ws.Range("N99").ClearContents()
$ws.Range("H132").Value = 4392.636
$ws.Range("I132").Value = 2685.5264
$ws.Range("K132").Value = 8056.5792
$ws.Range("M132").Value = -5526.5792

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H22").Value = 2540.889
$ws.Range("I22").Value = 2075
$ws.Range("J22").Value = 2913.6
$ws.Range("K22").Value = 2075
$ws.Range("L22").Value = 2913.6
$ws.Range("M22").Value = -1780
$ws.Range("N22").Value = -3503.6
$ws.Range("H27").Value = 2540.889
$ws.Range("I27").Value = 2075
$ws.Range("J27").Value = 2913.6
$ws.Range("K27").Value = 2075
$ws.Range("L27").Value = 2913.6
$ws.Range("M27").Value = -1968
$ws.Range("N27").Value = -3127.6
$ws.Range("H68").Value = 2181.182
$ws.Range("I68").Value = 2199.4285
$ws.Range("J68").Value = 2149.25
$ws.Range("K68").Value = 2199.4285
$ws.Range("L68").Value = 2149.25
$ws.Range("M68").Value = -1450.4285
$ws.Range("N68").Value = -3647.25
$ws.Range("H71").Value = 2181.182
$ws.Range("I71").Value = 2199.4285
$ws.Range("J71").Value = 2149.25
$ws.Range("K71").Value = 10997.1425
$ws.Range("L71").Value = 10746.25
$ws.Range("M71").Value = -7253.1425
$ws.Range("N71").Value = -18234.25
$ws.Range("H93").Value = 1729.561
$ws.Range("I93").Value = 1287.8077
$ws.Range("J93").Value = 2495.2666
$ws.Range("K93").Value = 1287.8077
$ws.Range("L93").Value = 2495.2666
$ws.Range("M93").Value = -39.80770000000007
$ws.Range("N93").Value = -4991.2666
$ws.Range("H100").Value = 3765.889
$ws.Range("I100").Value = 3301.4614
$ws.Range("K100").Value = 3301.4614
$ws.Range("M100").Value = -2760.4614
$ws.Range("H136").Value = 2039.8572
$ws.Range("I136").Value = 1893.56
$ws.Range("J136").Value = 3259
$ws.Range("K136").Value = 5680.68
$ws.Range("L136").Value = 9777
$ws.Range("M136").Value = -3130.68
$ws.Range("N136").Value = -14877

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H62").Value = 5959.8667
$ws.Range("I62").Value = 3780
$ws.Range("J62").Value = 7049.8
$ws.Range("K62").Value = 3780
$ws.Range("L62").Value = 7049.8
$ws.Range("M62").Value = -3156
$ws.Range("N62").Value = -8297.799999999999
$ws.Range("H65").Value = 5959.8667
$ws.Range("I65").Value = 3780
$ws.Range("J65").Value = 7049.8
$ws.Range("K65").Value = 18900
$ws.Range("L65").Value = 35249
$ws.Range("M65").Value = -15780
$ws.Range("N65").Value = -41489
$ws.Range("H81").Value = 1085.375
$ws.Range("I81").Value = 1149.8334
$ws.Range("K81").Value = 2299.6668
$ws.Range("M81").Value = -1238.6668
$ws.Range("H84").Value = 1085.375
$ws.Range("I84").Value = 1149.8334
$ws.Range("K84").Value = 11498.334
$ws.Range("M84").Value = -6194.333999999999
$ws.Range("H122").Value = 37075490
$ws.Range("I122").Value = 45501376
$ws.Range("K122").Value = 136504128
$ws.Range("M122").Value = -136501678
$ws.Range("H132").Value = 1866.3334
$ws.Range("I132").Value = 1819.4667
$ws.Range("J132").Value = 2335
$ws.Range("K132").Value = 5458.4001
$ws.Range("L132").Value = 7005
$ws.Range("M132").Value = -2928.4001
$ws.Range("N132").Value = -12065
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
